# Update peptide-based analysis (version 2)
# Replace the sumIntensity_N values in column A with the MSqRob
# peptide_abundance_study_variable.N. values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "peptide_abundance_study_variable.1."
$ws.Range("A3").Value = "peptide_abundance_study_variable.2."
$ws.Range("A4").Value = "peptide_abundance_study_variable.3."

# Move the active selection (matches the state captured in the saved file).
$ws.Range("G16").Select()
